# Regenerate the "K" column (column G) values from newly-computed strike counts.
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1; 3 = 2; 4 = 1; 5 = 0; 6 = 1; 7 = 0; 8 = 1; 9 = 1; 10 = 2;
    11 = 1; 12 = 0; 13 = 0; 14 = 0; 16 = 0; 17 = 2; 18 = 0;
    20 = 1; 21 = 2; 22 = 2; 23 = 2; 24 = 3; 25 = 1; 26 = 4; 27 = 2; 28 = 1; 29 = 0;
    30 = 4; 31 = 4; 32 = 1; 33 = 1; 34 = 1; 35 = 0; 36 = 0; 37 = 1; 38 = 0; 39 = 1;
    40 = 1; 41 = 1; 42 = 1; 43 = 2; 44 = 0; 45 = 1; 46 = 2; 47 = 0; 48 = 4; 49 = 3;
    50 = 0; 51 = 2; 52 = 1; 53 = 3; 54 = 1; 55 = 1; 56 = 1; 57 = 3; 58 = 1; 59 = 2;
    60 = 2; 61 = 1; 63 = 2; 64 = 0; 65 = 2; 66 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
